# Apply the cryptos.xlsx update: refreshed prices/volumes and a couple of
# rank swaps (rows 15/16, and a new BitDAO row inserted at 22 pushing
# Dai/Cosmos/WrappedBTC/... down by one row) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.044.74'
$ws.Range("E2").Value = '  -3.56%  '
$ws.Range("D3").Value = '1.602.04'
$ws.Range("E3").Value = '  -2.96%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''1.003'
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").Value = '''301.32'
$ws.Range("D7").Value = '''0.3781'
$ws.Range("E7").Value = '  -2.62%  '
$ws.Range("D8").Value = '''0.3631'
$ws.Range("E8").Value = '  -5.25%  '
$ws.Range("D9").Value = '''49.81'
$ws.Range("E9").Value = '  -2.77%  '
$ws.Range("D10").Value = '''1.260'
$ws.Range("E10").Value = '  -5.93%  '
$ws.Range("D11").Value = '''1.004'
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("D12").Value = '''0.08111'
$ws.Range("E12").Value = '  -3.87%  '
$ws.Range("E13").Value = '  -4.61%  '
$ws.Range("D14").Value = '''6.586'
$ws.Range("E14").Value = '  -6.00%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '''7.395'
$ws.Range("E15").Value = '  -7.73%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.00001254'
$ws.Range("E16").Value = '  -4.73%  '
$ws.Range("D17").Value = '1.596.81'
$ws.Range("E17").Value = '  -3.11%  '
$ws.Range("D18").Value = '''91.95'
$ws.Range("E18").Value = '  -2.23%  '
$ws.Range("D19").Value = '''0.06889'
$ws.Range("E19").Value = '  -1.27%  '
$ws.Range("D20").Value = '''18.21'
$ws.Range("E20").Value = '  -6.77%  '
$ws.Range("D21").Value = '''6.561'
$ws.Range("E21").Value = '  -5.58%  '
$ws.Range("B22").Value = 'BitDAO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range("D22").Value = '''0.5560'
$ws.Range("E22").Value = '  -5.33%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '''1.003'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '''13.11'
$ws.Range("E24").Value = '  -3.81%  '
$ws.Range("B25").Value = 'WrappedBTC'
$ws.Range("C25").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D25").Value = '23.006.03'
$ws.Range("E25").Value = '  -3.67%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '''2.369'
$ws.Range("E26").Value = '  -3.13%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '''2.787'
$ws.Range("E27").Value = '  -4.51%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '''21.08'
$ws.Range("E28").Value = '  -3.97%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '''150.50'
$ws.Range("E29").Value = '  -1.89%  '
$ws.Range("B30").Value = 'HuobiToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D30").Value = '''5.247'
$ws.Range("E30").Value = '  -2.64%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = '''133.18'
$ws.Range("E31").Value = '  -2.95%  '
$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").Value = '''2.324'
$ws.Range("E32").Value = '  -6.54%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''6.828'
$ws.Range("E33").Value = '  -11.63%  '
$ws.Range("B34").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C34").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D34").Value = '1.777.26'
$ws.Range("E34").Value = '  -3.57%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '''0.9625'
$ws.Range("E35").Value = '  -2.95%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '''0.07661'
$ws.Range("E36").Value = '  -6.00%  '
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").Value = '''10.40'
$ws.Range("E37").Value = '  -1.21%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '''6.296'
$ws.Range("E38").Value = '  -5.38%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.02718'
$ws.Range("E39").Value = '  -6.59%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '''0.2534'
$ws.Range("E40").Value = '  -5.38%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '''0.08854'
$ws.Range("E41").Value = '  -2.70%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''1.365'
$ws.Range("E42").Value = '  -4.08%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '''0.7041'
$ws.Range("E43").Value = '  -6.80%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '''12.54'
$ws.Range("E44").Value = '  -6.75%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''15.27'
$ws.Range("E45").Value = '  -8.30%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.6602'
$ws.Range("E46").Value = '  -4.75%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '''2.316'
$ws.Range("E47").Value = '  -5.22%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").Value = '''1.001'
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").Value = '''3.996'
$ws.Range("E49").Value = '  -2.43%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '''132.67'
$ws.Range("E50").Value = '  -0.83%  '
$ws.Range("D51").Value = '''0.07924'
